# Daily attendance processing - 2025-12-10 11:50:32
# Normalizes the "Recorded By" (column G) value ordering on the
# "Session Analysis Results" sheet so that specific known orderings
# are rewritten to their corrected order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of exact old "Recorded By" text -> new corrected text.
$map = @{
    "System, dnasr281@gmail.com"            = "dnasr281@gmail.com, System"
    "System, admin@admin.com"               = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com"   = "admin@admin.com, dnasr281@gmail.com"
    "backup@backdoor.com, System, system"   = "system, backup@backdoor.com, System"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
